$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ADIOS Data Model Version"
$ws.Range("B1").Value = "0.11.0"
$ws.Range("A2").Value = "Record Metadata"
$ws.Range("A3").Value = "    Name"
$ws.Range("A4").Value = "    Source ID"
$ws.Range("A5").Value = "    Alternate Names"
$ws.Range("A6").Value = "    Location"
$ws.Range("A7").Value = "    Reference"
$ws.Range("B7").Value = "year"
$ws.Range("C7").Value = "full reference"
$ws.Range("A8").Value = "    Sample Date"
$ws.Range("A9").Value = "    Product Type"
$ws.Range("A10").Value = "    API"
$ws.Range("A11").Value = "    Labels"
$ws.Range("A12").Value = "    Location Coordinates"
$ws.Range("A13").Value = "    Comments"
$ws.Range("A15").Value = "Subsample Metadata"
$ws.Range("A16").Value = "    Name"
$ws.Range("A17").Value = "    Short name"
$ws.Range("A18").Value = "    Sample ID"
$ws.Range("A19").Value = "    Description"
$ws.Range("A20").Value = "    Fraction evaporated"
$ws.Range("A21").Value = "    Boiling Point Range"
$ws.Range("A23").Value = "Physical Properties"
$ws.Range("A24").Value = "  Pour Point"
$ws.Range("A25").Value = "  Flash Point"
$ws.Range("A26").Value = "  Density"
$ws.Range("A27").Value = "    Density at temp"
$ws.Range("A28").Value = "    Density at temp"
$ws.Range("A29").Value = "    Density at temp"
$ws.Range("A31").Value = "Viscosity"
$ws.Range("A32").Value = "    Viscosity at temp"
$ws.Range("A33").Value = "    Viscosity at temp"
$ws.Range("A34").Value = "    Viscosity at temp"
$ws.Range("A36").Value = "Distillation Data"
$ws.Range("A37").Value = "  Type (mass fraction or volume fraction)"
$ws.Range("A38").Value = "    Method"
$ws.Range("A39").Value = "    Final Boiling point"
$ws.Range("A40").Value = "    Fraction Recovered"
$ws.Range("A41").Value = "  Distillation cuts"
$ws.Range("B41").Value = "Fraction"
$ws.Range("D41").Value = "Temp"
$ws.Range("E41").Value = "Temp Unit"
$ws.Range("A42").Value = "    cut 1"
$ws.Range("A43").Value = "    cut 2"
$ws.Range("A44").Value = "    cut 3"
$ws.Range("A46").Value = "  SARA Analysis"
$ws.Range("B46").Value = "Fraction"
$ws.Range("C46").Value = "Fraction Unit"
$ws.Range("A47").Value = "    Method"
$ws.Range("A48").Value = "    Saturates"
$ws.Range("A49").Value = "    Aromatics"
$ws.Range("A50").Value = "    Resins"
$ws.Range("A51").Value = "    Asphaltenes"

$boldRng = $ws.Range("A1,A2,A15,A23,A36,A46")
$boldRng.Font.Bold = $true

$altRng = $ws.Range("A47,A48,A49,A50,A51,A52,A54,A55,A57,A58,A60,A61")
$altRng.Font.Name = "Arial"
$altRng.Font.Size = 11
$altRng.Font.Color = 0

$ws.Columns("C").ColumnWidth = 11.75

$ws.Range("C46").Select()
